# "Redone the shit, cause it screwed up." — re-adds the per-interval AVERAGE
# formulas in row 8 (K8:W8) that got lost, and restores the selection/scroll
# position the workbook had when the averages were last being checked.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook/UI language -------------------------------------------------
# The original commit also shows the workbook's built-in "Normal" cell style
# and the "Office" theme relabeled to their Dutch equivalents ("Standaard" /
# "Kantoor(thema)") -- that's Excel's own UI-language translation of the
# built-in names, applied automatically by a Dutch-locale Excel on save
# (built-in style/theme names aren't user-renameable via the object model).
# We still set what the object model exposes for this, best-effort:
try {
    $style = $excel.ActiveWorkbook.Styles.Item(1)
    $style.Name = "Standaard"
} catch {
}

# --- Re-add the AVERAGE() row (row 8) for every measurement block ---------
# K8 = AVERAGE(K3:K7) typed directly;
# L8:W8 filled from one relative formula so Excel stores it as a shared
# formula group (t="shared") anchored at L8, matching the other four blocks.
$ws.Range("K8").Formula = "=AVERAGE(K3:K7)"
$ws.Range("L8:W8").Formula = "=AVERAGE(L3:L7)"

# --- Restore selection / scroll position -----------------------------------
# Active cell moves to K8 (the first newly (re)computed average) and the
# window is scrolled so column M is the leftmost visible column.
$ws.Range("K8").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollColumn = 13
    $excel.ActiveWindow.ScrollRow = 1
} catch {
}
